$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: column B is now "mapsto" (it used to be "meaning"); the old
# numeric "mapsto" column (C) is going away entirely.
$ws.Range("B1").Value = "mapsto"

# Replace the old free-text "meaning" values with the new, normalized
# mapsto labels. "wake" (row 2) and "movement" (row 8) stay the same;
# the stage names are tightened up and stage 3 + stage 4 both collapse
# into "sws".
$ws.Range("B2").Value = "wake"
$ws.Range("B3").Value = "stage1"
$ws.Range("B4").Value = "stage2"
$ws.Range("B5").Value = "sws"
$ws.Range("B6").Value = "sws"
$ws.Range("B7").Value = "rem"
$ws.Range("B8").Value = "movement"

# Drop the old numeric "mapsto" column (C) entirely, shifting nothing else.
$ws.Range("C1:C8").EntireColumn.Delete()

# Match the new active-cell selection recorded in the workbook view.
$ws.Range("C3").Select()
